$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Manchester tribunal details ---
# Address (row 3) - text update only
$ws.Range("B3").Value = "Manchester Employment Tribunal, Alexandra House, 14-22 The Parsonage, Manchester, M3 2JA"

# Telephone / Fax / DX (rows 4-6) - were bare numbers, now become text values
$ws.Range("B4").Value = "0161 833 6100"
$ws.Range("B5").Value = "0870 739 4433"
$ws.Range("B6").Value = "DX 743570"

# Email (row 7)
$ws.Range("B7").Value = "Manchesteret@justice.gov.uk"

# --- Glasgow tribunal details ---
# Address (row 8) - text update only
$ws.Range("B8").Value = "Eagle Building, 215 Bothwell Street, Glasgow, G2 7TS"

# Telephone / Fax / DX (rows 9-11) - were bare numbers, now become text values
$ws.Range("B9").Value = "0141 204 0730"
$ws.Range("B10").Value = "01264 785 177"
$ws.Range("B11").Value = "DX 7435701"

# Email (row 12)
$ws.Range("B12").Value = "glasgowet@justice.gov.uk"

# --- Hyperlinks: rebuild so the display text + mailto targets match the new emails ---
# (Individual Hyperlink.TextToDisplay/.Address writes or Item(n).Delete() do not mutate
# the existing entries in-place on this host, so the collection is cleared and the two
# mailto links are re-added pointing at the updated addresses.)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:Manchesteret@justice.gov.uk", "", "", "Manchesteret@justice.gov.uk")
$ws.Hyperlinks.Add($ws.Range("B12"), "mailto:glasgowet@justice.gov.uk", "", "", "glasgowet@justice.gov.uk")
